$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pipeline(steps=[('scaler', RobustScaler()),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 MLPClassifier(alpha=1, hidden_layer_sizes=(10,),`n                               learning_rate_init=0.01, max_iter=1000,`n                               random_state=42, solver='sgd'))])"
$ws.Range("B2").Value = 0.6761904761904762
$ws.Range("C2").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': RobustScaler(), 'model__solver': 'sgd', 'model__learning_rate_init': 0.01, 'model__hidden_layer_sizes': (10,), 'model__alpha': 1, 'model__activation': 'relu'}"
$ws.Range("D2").Value = 0.1666666666666667
$ws.Range("E2").Value = "[1 1 0 0 1 0 0 0 0 1 0 1]"
$ws.Range("F2").Value = "[0 0 1 1 0 1 1 0 1 0 1 1]"
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.7221207087486158
$ws.Range("I2").Value = 0.03945522122388282
$ws.Range("J2").Value = 0.5547065337763012
$ws.Range("K2").Value = 0.07818516010006489
$ws.Rows(2).AutoFit()

$ws.Range("A3").Value = "Pipeline(steps=[('scaler', StandardScaler()),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 MLPClassifier(activation='tanh', alpha=0.01,`n                               hidden_layer_sizes=(10,),`n                               learning_rate_init=0.0001, max_iter=1000,`n                               random_state=42, solver='sgd'))])"
$ws.Range("B3").Value = 0.6380952380952382
$ws.Range("C3").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': StandardScaler(), 'model__solver': 'sgd', 'model__learning_rate_init': 0.0001, 'model__hidden_layer_sizes': (10,), 'model__alpha': 0.01, 'model__activation': 'tanh'}"
$ws.Range("D3").Value = 0.7999999999999999
$ws.Range("E3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0]"
$ws.Range("F3").Value = "[1 1 0 1 0 1 0 1 1 1 1 0]"
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.7252725186460126
$ws.Range("I3").Value = 0.04069783648809968
$ws.Range("J3").Value = 0.5287435456110154
$ws.Range("K3").Value = 0.07713770913083769
$ws.Rows(3).AutoFit()

$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None), ('selector', None),`n                ('model',`n                 MLPClassifier(activation='tanh', alpha=1e-05,`n                               hidden_layer_sizes=(5, 10, 5),`n                               learning_rate_init=0.0001, max_iter=1000,`n                               random_state=42))])"
$ws.Range("B4").Value = 0.6
$ws.Range("C4").Value = "{'selector': None, 'scaler': None, 'model__solver': 'adam', 'model__learning_rate_init': 0.0001, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 1e-05, 'model__activation': 'tanh'}"
$ws.Range("D4").Value = 0.8
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 1 1 1 1 1 1 1 1 1 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.7501904761904762
$ws.Range("I4").Value = 0.04546878142271138
$ws.Range("J4").Value = 0.528126984126984
$ws.Range("K4").Value = 0.07809923984026701
$ws.Rows(4).AutoFit()

$ws.Range("A5").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),`n                ('model',`n                 MLPClassifier(alpha=0.01, hidden_layer_sizes=(5, 10, 5),`n                               learning_rate_init=0.0001, max_iter=1000,`n                               random_state=42, solver='lbfgs'))])"
$ws.Range("B5").Value = 0.6380952380952382
$ws.Range("C5").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__solver': 'lbfgs', 'model__learning_rate_init': 0.0001, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 0.01, 'model__activation': 'relu'}"
$ws.Range("D5").Value = 0.8235294117647058
$ws.Range("E5").Value = "[1 1 0 0 0 0 1 0 1 1 1 1]"
$ws.Range("F5").Value = "[1 1 1 1 0 1 1 0 1 1 1 1]"
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.739342403628118
$ws.Range("I5").Value = 0.04613195158323608
$ws.Range("J5").Value = 0.5060090702947846
$ws.Range("K5").Value = 0.08697345844275788
$ws.Rows(5).AutoFit()

$ws.Range("A6").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',`n                                                     random_state=42))),`n                ('model',`n                 MLPClassifier(activation='tanh', alpha=0.01,`n                               hidden_layer_sizes=(5, 10, 5),`n                               learning_rate_init=0.0001, max_iter=1000,`n                               random_state=42))])"
$ws.Range("B6").Value = 0.6285714285714287
$ws.Range("C6").Value = "{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__solver': 'adam', 'model__learning_rate_init': 0.0001, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 0.01, 'model__activation': 'tanh'}"
$ws.Range("D6").Value = 0.4615384615384615
$ws.Range("E6").Value = "[1 1 1 1 0 0 0 0 1 1 0 0]"
$ws.Range("F6").Value = "[1 1 0 0 0 1 1 1 1 0 1 0]"
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.7285420340975897
$ws.Range("I6").Value = 0.04375204601354599
$ws.Range("J6").Value = 0.5490887713109934
$ws.Range("K6").Value = 0.07445945632576788
$ws.Rows(6).AutoFit()

